$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "256.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.07%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.34%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.682"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.88%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05940"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.24%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.618"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.71%"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.86%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9248"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.42%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1379"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.67%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04406"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "16.27%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07002"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.56%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03059"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.50%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09097"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001531"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.78%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006189"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.72%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.470"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.77%"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.161"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.30%"

$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.199"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.91%"

$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.01019"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1,590.57%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3028"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.71%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.30%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.856"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.30%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04246"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.48%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001219"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.38%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004767"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "11.20%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.06%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2.10%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03771"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.40%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006282"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.52%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1094"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.68%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002201"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.06%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01410"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "26.46%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005314"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.36%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.06%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.04402"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-50.26%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2415"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "9,816.50%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.06%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.06%"
